# error solve ifrs list
# Replace the per-company financial figures in rows 2-6 with corrected values,
# and remove the (now invalid) data for rows 7-9, keeping only the
# identifier columns (A: index, B: ticker/code, C: company name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New figures for rows 2-6, keyed by column letter -> value.
$rowUpdates = @{
    2 = @{
        D = 13237; E = 163; F = 163; G = -71; H = -80; I = -54; J = -25
        K = 6557; L = 4877; M = 1679; N = 1579; O = 100; P = 456; Q = 277
        R = -93; S = -240; T = 111; U = 167; V = 4129
        W = 1.23; X = -0.6; Y = -3.36; Z = -1.22
        AA = 290.46; AB = 266.04; AC = -60; AD = -14.73; AE = 2053
        AF = 0.43; AG = 0; AH = 0; AI = 0
        AJ = 91140499
    }
    3 = @{
        D = 10167; E = -320; F = -320; G = -526; H = -509; I = -447; J = -61
        K = 5391; L = 3900; M = 1491; N = 1427; O = 64; P = 456; Q = 1048
        R = -66; S = -653; T = 275; U = 773; V = 3483
        W = -3.15; X = -5; Y = -29.76; Z = -8.52
        AA = 261.6; AB = 164.99; AC = -491; AD = -2.22; AE = 1856
        AF = 0.59; AG = 0; AH = 0; AI = 0
        AJ = 91140499
    }
    4 = @{
        D = 8177; E = 302; F = 302; G = 118; H = 120; I = 92; J = 28
        K = 5660; L = 3867; M = 1792; N = 1473; O = 320; P = 456; Q = 429
        R = -73; S = -623; T = 89; U = 340; V = 3178
        W = 3.7; X = 1.47; Y = 6.34; Z = 2.17
        AA = 215.76; AB = 177.55; AC = 101; AD = 9.880000000000001; AE = 1915
        AF = 0.52; AG = 0; AH = 0; AI = 0
        AJ = 91140499
    }
    5 = @{
        D = 9833; E = 406; F = 406; G = 162; H = 143; I = 127; J = 15
        K = 5634; L = 3689; M = 1945; N = 1610; O = 335; P = 456; Q = 108
        R = -36; S = -63; T = 23; U = 84; V = 3094
        W = 4.13; X = 1.45; Y = 8.26; Z = 2.53
        AA = 189.68; AB = 206.5; AC = 140; AD = 7.55; AE = 2093
        AF = 0.5; AG = 25; AH = 2.37; AI = 15.1
        AJ = 91140499
    }
    6 = @{
        D = 9513; E = 173; F = 173; G = 135; H = 124; I = 98
        K = 5700; L = 3542; M = 2157; N = 1835; P = 456; Q = 288
        R = -90; S = -156; T = 95; U = 193; V = 2906
        W = 1.82; X = 1.31; Y = 5.7; Z = 2.19
        AA = 164.19; AB = 232.77; AC = 108; AD = 8.58; AE = 2385
        AF = 0.39; AG = 30; AH = 3.24; AI = 23.49
        AJ = 91140499
    }
}

# Columns in sheet order, used only to keep writes in a deterministic order.
$colOrder = @('D','E','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T','U','V', `
              'W','X','Y','Z','AA','AB','AC','AD','AE','AF','AG','AH','AI','AJ')

foreach ($rowNum in ($rowUpdates.Keys | Sort-Object)) {
    $cols = $rowUpdates[$rowNum]
    foreach ($col in $colOrder) {
        if ($cols.ContainsKey($col)) {
            $ws.Range("$col$rowNum").Value = $cols[$col]
        }
    }
}

# Rows 7-9 no longer have usable figures; clear D:AJ but keep A (index),
# B (ticker) and C (company name) intact.
foreach ($rowNum in 7..9) {
    $ws.Range("D${rowNum}:AJ${rowNum}").ClearContents()
}
